# Add team record columns (Wins, Losses, Ties) to the sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new columns AD, AE, AF
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy header formatting (bold, border, centered) from an existing header cell
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Fill in team record data for every data row (rows 2-51)
for ($r = 2; $r -le 51; $r++) {
    $ws.Cells.Item($r, 30).Value = 82  # AD - Wins
    $ws.Cells.Item($r, 31).Value = 80  # AE - Losses
    $ws.Cells.Item($r, 32).Value = 0   # AF - Ties
}
